# Sheet4 ("存款" / deposit) restructuring: add bank/deposit_type/currency
# columns and the common property/category/date/legislator/index trailer
# columns (G-M), matching the layout used on the other property sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1) ---
$header = New-Object 'object[,]' 1,12
$header[0,0] = "bank"
$header[0,1] = "deposit_type"
$header[0,2] = "currency"
$header[0,3] = "owner"
$header[0,4] = "total"
$header[0,5] = "property_category"
$header[0,6] = "category"
$header[0,7] = "date"
$header[0,8] = "legislator_name"
$header[0,9] = "legislator_id"
$header[0,10] = "source_file"
$header[0,11] = "index"
$ws.Range("B1:M1").Value = $header

# --- Data rows 2-18 (columns A-M) ---
$data = New-Object 'object[,]' 17,13
$data[0,0] = 49
$data[0,1] = "中華郵政股份有限公司"
$data[0,2] = "活期存款"
$data[0,3] = "新臺幣"
$data[0,4] = "江啟臣"
$data[0,5] = 4745947
$data[0,6] = "deposit"
$data[0,7] = "normal"
$data[0,8] = "2012-04-20"
$data[0,9] = "江啟臣"
$data[0,10] = 1731
$data[0,11] = "tmpe0681"
$data[0,12] = 49
$data[1,0] = 50
$data[1,1] = "中國信託商業銀行"
$data[1,2] = "活期存款"
$data[1,3] = "新臺幣"
$data[1,4] = "江啟臣"
$data[1,5] = 2616299
$data[1,6] = "deposit"
$data[1,7] = "normal"
$data[1,8] = "2012-04-20"
$data[1,9] = "江啟臣"
$data[1,10] = 1731
$data[1,11] = "tmpe0681"
$data[1,12] = 50
$data[2,0] = 51
$data[2,1] = "國泰世華商業銀行"
$data[2,2] = "活期儲蓄存款"
$data[2,3] = "新臺幣"
$data[2,4] = "江啟臣"
$data[2,5] = 668336
$data[2,6] = "deposit"
$data[2,7] = "normal"
$data[2,8] = "2012-04-20"
$data[2,9] = "江啟臣"
$data[2,10] = 1731
$data[2,11] = "tmpe0681"
$data[2,12] = 51
$data[3,0] = 52
$data[3,1] = "國泰世華商業銀行"
$data[3,2] = "活期儲蓄存款"
$data[3,3] = "新臺幣"
$data[3,4] = "江啟臣"
$data[3,5] = 1
$data[3,6] = "deposit"
$data[3,7] = "normal"
$data[3,8] = "2012-04-20"
$data[3,9] = "江啟臣"
$data[3,10] = 1731
$data[3,11] = "tmpe0681"
$data[3,12] = 52
$data[4,0] = 53
$data[4,1] = "臺灣銀行"
$data[4,2] = "綜合存款"
$data[4,3] = "新臺幣"
$data[4,4] = "江啟臣"
$data[4,5] = 75103
$data[4,6] = "deposit"
$data[4,7] = "normal"
$data[4,8] = "2012-04-20"
$data[4,9] = "江啟臣"
$data[4,10] = 1731
$data[4,11] = "tmpe0681"
$data[4,12] = 53
$data[5,0] = 54
$data[5,1] = "匯豐(台灣)商業銀行"
$data[5,2] = "活期存款"
$data[5,3] = "新臺幣"
$data[5,4] = "江啟臣"
$data[5,5] = 500075
$data[5,6] = "deposit"
$data[5,7] = "normal"
$data[5,8] = "2012-04-20"
$data[5,9] = "江啟臣"
$data[5,10] = 1731
$data[5,11] = "tmpe0681"
$data[5,12] = 54
$data[6,0] = 55
$data[6,1] = "臺灣銀行"
$data[6,2] = "活期儲蓄存款"
$data[6,3] = "新臺幣"
$data[6,4] = "劉姿伶"
$data[6,5] = 15710
$data[6,6] = "deposit"
$data[6,7] = "normal"
$data[6,8] = "2012-04-20"
$data[6,9] = "江啟臣"
$data[6,10] = 1731
$data[6,11] = "tmpe0681"
$data[6,12] = 55
$data[7,0] = 56
$data[7,1] = "台新國際商業銀行"
$data[7,2] = "活期儲蓄存款"
$data[7,3] = "新臺幣"
$data[7,4] = "劉姿伶"
$data[7,5] = 12270
$data[7,6] = "deposit"
$data[7,7] = "normal"
$data[7,8] = "2012-04-20"
$data[7,9] = "江啟臣"
$data[7,10] = 1731
$data[7,11] = "tmpe0681"
$data[7,12] = 56
$data[8,0] = 57
$data[8,1] = "彰化商業銀行"
$data[8,2] = "活期存款"
$data[8,3] = "新臺幣"
$data[8,4] = "劉姿伶"
$data[8,5] = 153468
$data[8,6] = "deposit"
$data[8,7] = "normal"
$data[8,8] = "2012-04-20"
$data[8,9] = "江啟臣"
$data[8,10] = 1731
$data[8,11] = "tmpe0681"
$data[8,12] = 57
$data[9,0] = 58
$data[9,1] = "匯豐(台灣)商業銀行"
$data[9,2] = "活期儲蓄存款"
$data[9,3] = "新臺幣"
$data[9,4] = "劉姿伶"
$data[9,5] = 188120
$data[9,6] = "deposit"
$data[9,7] = "normal"
$data[9,8] = "2012-04-20"
$data[9,9] = "江啟臣"
$data[9,10] = 1731
$data[9,11] = "tmpe0681"
$data[9,12] = 58
$data[10,0] = 59
$data[10,1] = "合作金庫商業銀行"
$data[10,2] = "活期儲蓄存款"
$data[10,3] = "新臺幣"
$data[10,4] = "劉姿伶"
$data[10,5] = 222237
$data[10,6] = "deposit"
$data[10,7] = "normal"
$data[10,8] = "2012-04-20"
$data[10,9] = "江啟臣"
$data[10,10] = 1731
$data[10,11] = "tmpe0681"
$data[10,12] = 59
$data[11,0] = 60
$data[11,1] = "合作金庫商業銀行"
$data[11,2] = "活期儲蓄存款"
$data[11,3] = "新臺幣"
$data[11,4] = "劉姿伶"
$data[11,5] = 695410
$data[11,6] = "deposit"
$data[11,7] = "normal"
$data[11,8] = "2012-04-20"
$data[11,9] = "江啟臣"
$data[11,10] = 1731
$data[11,11] = "tmpe0681"
$data[11,12] = 60
$data[12,0] = 61
$data[12,1] = "國泰世華商業銀行"
$data[12,2] = "活期儲蓄存款"
$data[12,3] = "新臺幣"
$data[12,4] = "劉姿伶"
$data[12,5] = 212669
$data[12,6] = "deposit"
$data[12,7] = "normal"
$data[12,8] = "2012-04-20"
$data[12,9] = "江啟臣"
$data[12,10] = 1731
$data[12,11] = "tmpe0681"
$data[12,12] = 61
$data[13,0] = 62
$data[13,1] = "國泰世華商業銀行"
$data[13,2] = "活期儲蓄存款"
$data[13,3] = "新臺幣"
$data[13,4] = "劉姿伶"
$data[13,5] = 773
$data[13,6] = "deposit"
$data[13,7] = "normal"
$data[13,8] = "2012-04-20"
$data[13,9] = "江啟臣"
$data[13,10] = 1731
$data[13,11] = "tmpe0681"
$data[13,12] = 62
$data[14,0] = 63
$data[14,1] = "國泰世華商業銀行"
$data[14,2] = "活期儲蓄存款"
$data[14,3] = "美金"
$data[14,4] = "劉姿伶"
$data[14,5] = 61395
$data[14,6] = "deposit"
$data[14,7] = "normal"
$data[14,8] = "2012-04-20"
$data[14,9] = "江啟臣"
$data[14,10] = 1731
$data[14,11] = "tmpe0681"
$data[14,12] = 63
$data[15,0] = 64
$data[15,1] = "台北富邦商業銀行"
$data[15,2] = "活期儲蓄存款"
$data[15,3] = "新臺幣"
$data[15,4] = "劉姿伶"
$data[15,5] = 2569
$data[15,6] = "deposit"
$data[15,7] = "normal"
$data[15,8] = "2012-04-20"
$data[15,9] = "江啟臣"
$data[15,10] = 1731
$data[15,11] = "tmpe0681"
$data[15,12] = 64
$data[16,0] = 65
$data[16,1] = "中華郵政股份有限公司"
$data[16,2] = "活期儲蓄存款"
$data[16,3] = "新臺幣"
$data[16,4] = "劉姿伶"
$data[16,5] = 360747
$data[16,6] = "deposit"
$data[16,7] = "normal"
$data[16,8] = "2012-04-20"
$data[16,9] = "江啟臣"
$data[16,10] = 1731
$data[16,11] = "tmpe0681"
$data[16,12] = 65

# The "date" column (I) holds a text value ("2012-04-20") elsewhere in this
# workbook, not a real Excel date -- force text interpretation so COM does
# not auto-convert the string into a date serial number.
$ws.Range("I2:I18").NumberFormat = "@"
$ws.Range("A2:M18").Value = $data
$ws.Range("I2:I18").Style = "Normal"

